$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.801.52'
$ws.Range("E2").Value = '  +2.75%  '
$ws.Range("D3").Value = '1.866.42'
$ws.Range("E3").Value = '  +2.42%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.039'
$ws.Range("E4").Value = '  +2.75%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.76'
$ws.Range("E5").Value = '  +3.17%  '
$ws.Range("E6").Value = '  +2.44%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4422'
$ws.Range("E7").Value = '  +2.48%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3802'
$ws.Range("E8").Value = '  +2.62%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07476'
$ws.Range("E9").Value = '  +2.61%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8863'
$ws.Range("E10").Value = '  +1.56%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.80'
$ws.Range("E11").Value = '  +2.04%  '
$ws.Range("D12").Value = '1.886.75'
$ws.Range("E12").Value = '  -12.59%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.564'
$ws.Range("E13").Value = '  +2.58%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.764'
$ws.Range("E14").Value = '  +1.53%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07240'
$ws.Range("E15").Value = '  +3.33%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '83.85'
$ws.Range("E16").Value = '  +3.15%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.041'
$ws.Range("E17").Value = '  +2.88%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009160'
$ws.Range("E18").Value = '  +2.42%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.036'
$ws.Range("E19").Value = '  +2.53%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.58'
$ws.Range("E20").Value = '  +1.41%  '
$ws.Range("D21").Value = '27.818.23'
$ws.Range("E21").Value = '  +2.58%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.321'
$ws.Range("E22").Value = '  +1.92%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.39'
$ws.Range("E23").Value = '  +3.22%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.967'
$ws.Range("E24").Value = '  +3.88%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.57'
$ws.Range("E25").Value = '  +2.36%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.92'
$ws.Range("E26").Value = '  +2.49%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.989'
$ws.Range("E27").Value = '  +3.37%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.330'
$ws.Range("E28").Value = '  +1.41%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '117.85'
$ws.Range("E29").Value = '  +2.33%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09116'
$ws.Range("E30").Value = '  +1.36%  '
$ws.Range("B31").Value = 'ARBITRUM'
$ws.Range("C31").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.223'
$ws.Range("E31").Value = '  +3.49%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7794'
$ws.Range("E32").Value = '  +3.89%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.059'
$ws.Range("E33").Value = '  +8.26%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.593'
$ws.Range("E34").Value = '  +3.33%  '
$ws.Range("E35").Value = '  +2.66%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.165'
$ws.Range("E36").Value = '  +3.42%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01999'
$ws.Range("E37").Value = '  +3.44%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05369'
$ws.Range("E38").Value = '  +2.15%  '
$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5209'
$ws.Range("E39").Value = '  +1.43%  '
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.847'
$ws.Range("E40").Value = '  +3.27%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1699'
$ws.Range("E41").Value = '  +2.55%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.927'
$ws.Range("E42").Value = '  +6.48%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.727'
$ws.Range("E43").Value = '  +4.35%  '
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '109.84'
$ws.Range("E44").Value = '  +2.18%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.69'
$ws.Range("E45").Value = '  +2.19%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.728'
$ws.Range("E46").Value = '  +4.14%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4723'
$ws.Range("E47").Value = '  +2.60%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06452'
$ws.Range("E48").Value = '  +3.35%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.890'
$ws.Range("E49").Value = '  +3.25%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '40.06'
$ws.Range("E50").Value = '  +4.54%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '64.68'
$ws.Range("E51").Value = '  +1.27%  '
